$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1300.4736
$ws.Range("J19").Value = 1206
$ws.Range("L19").Value = 1206
$ws.Range("N19").Value = -1556
# Row 32
$ws.Range("H32").Value = 1936.875
$ws.Range("I32").Value = 1833.3334
$ws.Range("J32").Value = 1999
$ws.Range("K32").Value = 1833.3334
$ws.Range("L32").Value = 1999
$ws.Range("M32").Value = -1507.3334
$ws.Range("N32").Value = -2651
# Row 53
$ws.Range("H53").Value = 352.4
$ws.Range("I53").Value = 117.416664
$ws.Range("J53").Value = 704.875
$ws.Range("K53").Value = 117.416664
$ws.Range("L53").Value = 704.875
$ws.Range("M53").Value = 519.583336
$ws.Range("N53").Value = -1978.875
# Row 62
$ws.Range("H62").Value = 4054.4
$ws.Range("I62").Value = 3499
$ws.Range("J62").Value = 4887.5
$ws.Range("K62").Value = 3499
$ws.Range("L62").Value = 4887.5
$ws.Range("M62").Value = -2875
$ws.Range("N62").Value = -6135.5
# Row 65
$ws.Range("H65").Value = 4054.4
$ws.Range("I65").Value = 3499
$ws.Range("J65").Value = 4887.5
$ws.Range("K65").Value = 17495
$ws.Range("L65").Value = 24437.5
$ws.Range("M65").Value = -14375
$ws.Range("N65").Value = -30677.5
# Row 80
$ws.Range("H80").Value = 1794.6
$ws.Range("I80").Value = 649.5
$ws.Range("J80").Value = 2558
$ws.Range("K80").Value = 1948.5
$ws.Range("L80").Value = 7674
$ws.Range("M80").Value = -950.5
$ws.Range("N80").Value = -9670
# Row 83
$ws.Range("H83").Value = 1794.6
$ws.Range("I83").Value = 649.5
$ws.Range("J83").Value = 2558
$ws.Range("K83").Value = 5845.5
$ws.Range("L83").Value = 23022
$ws.Range("M83").Value = -853.5
$ws.Range("N83").Value = -33006
# Row 86
$ws.Range("H86").Value = 7799.375
$ws.Range("J86").Value = 9331.666999999999
$ws.Range("L86").Value = 9331.666999999999
$ws.Range("N86").Value = -11577.667
# Row 88
$ws.Range("H88").Value = 2770.6667
$ws.Range("J88").Value = 2659.0908
$ws.Range("L88").Value = 2659.0908
$ws.Range("N88").Value = -3471.0908
# Row 89
$ws.Range("H89").Value = 7799.375
$ws.Range("J89").Value = 9331.666999999999
$ws.Range("L89").Value = 46658.335
$ws.Range("N89").Value = -57890.335
# Row 91
$ws.Range("H91").Value = 2770.6667
$ws.Range("J91").Value = 2659.0908
$ws.Range("L91").Value = 2659.0908
$ws.Range("N91").Value = -5467.0908
# Row 113
$ws.Range("H113").Value = 49999.5
$ws.Range("I113").Value = 49999.5
$ws.Range("K113").Value = 49999.5
$ws.Range("M113").Value = -46745.5
# Row 132
$ws.Range("H132").Value = 3349.7646
$ws.Range("I132").Value = 3261.6667
$ws.Range("K132").Value = 9785.000100000001
$ws.Range("M132").Value = -7255.000100000001
# Row 137
$ws.Range("H137").Value = 1638.7
$ws.Range("I137").Value = 1398.1538
$ws.Range("K137").Value = 4194.4614
$ws.Range("M137").Value = -1644.4614
# Row 138
$ws.Range("H138").Value = 2347.0833
$ws.Range("I138").Value = 2105.9092
$ws.Range("K138").Value = 6317.7276
$ws.Range("M138").Value = -1177.7276
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 30
$ws.Range("H30").Value = 15000
$ws.Range("J30").Value = 15000
$ws.Range("L30").Value = 15000
$ws.Range("N30").Value = -15300
# Row 32
$ws.Range("H32").Value = 8952.385
$ws.Range("I32").Value = 8113.143
$ws.Range("K32").Value = 8113.143
$ws.Range("M32").Value = -7826.143
# Row 50
$ws.Range("H50").Value = 2842.3635
$ws.Range("I50").Value = 1894.3334
$ws.Range("J50").Value = 3980
$ws.Range("K50").Value = 1894.3334
$ws.Range("L50").Value = 3980
$ws.Range("M50").Value = -1180.3334
$ws.Range("N50").Value = -5408
# Row 76
$ws.Range("H76").Value = 55000
$ws.Range("J76").Value = 55000
$ws.Range("L76").Value = 55000
$ws.Range("N76").Value = -55676
# Row 79
$ws.Range("H79").Value = 55000
$ws.Range("J79").Value = 55000
$ws.Range("L79").Value = 55000
$ws.Range("N79").Value = -57340
# Row 110
$ws.Range("H110").Value = 1547.4286
$ws.Range("I110").Value = 1547.4286
$ws.Range("K110").Value = 1547.4286
$ws.Range("M110").Value = 497.5714
# Row 121
$ws.Range("H121").Value = 98999
$ws.Range("J121").Value = 98999
$ws.Range("L121").Value = 98999
$ws.Range("N121").Value = -102493

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 13
$ws.Range("H13").Value = 50000
$ws.Range("J13").Value = 50000
$ws.Range("L13").Value = 50000
$ws.Range("N13").Value = -50336
# Row 64
$ws.Range("H64").Value = 933.375
$ws.Range("I64").Value = 897.25
$ws.Range("K64").Value = 897.25
$ws.Range("M64").Value = -672.25
# Row 67
$ws.Range("H67").Value = 933.375
$ws.Range("I67").Value = 897.25
$ws.Range("K67").Value = 897.25
$ws.Range("M67").Value = -117.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 128
$ws.Range("H128").Value = 277759
$ws.Range("I128").Value = 277759
$ws.Range("K128").Value = 833277
$ws.Range("M128").Value = -828297
# Row 131
$ws.Range("H131").Value = 1119.2
$ws.Range("J131").Value = 1141
$ws.Range("L131").Value = 3423
$ws.Range("N131").Value = -13503
# Row 140
$ws.Range("H140").Value = 1500
$ws.Range("I140").Value = 1500
$ws.Range("K140").Value = 4500
$ws.Range("M140").Value = 680

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 59
$ws.Range("H59").Value = 39998
$ws.Range("J59").Value = 39998
$ws.Range("L59").Value = 39998
$ws.Range("N59").Value = -41164
# Row 70
$ws.Range("H70").Value = 11299.8
$ws.Range("I70").Value = 11299.8
$ws.Range("K70").Value = 11299.8
$ws.Range("M70").Value = -11029.8
# Row 73
$ws.Range("H73").Value = 11299.8
$ws.Range("I73").Value = 11299.8
$ws.Range("K73").Value = 11299.8
$ws.Range("M73").Value = -10363.8
# Row 80
$ws.Range("H80").Value = 2888.7778
$ws.Range("I80").Value = 2758.1667
$ws.Range("K80").Value = 2758.1667
$ws.Range("M80").Value = -1760.1667
# Row 83
$ws.Range("H83").Value = 2888.7778
$ws.Range("I83").Value = 2758.1667
$ws.Range("K83").Value = 13790.8335
$ws.Range("M83").Value = -8798.833500000001
# Row 122
$ws.Range("H122").Value = 2699.6667
$ws.Range("I122").Value = 3050
$ws.Range("K122").Value = 9150
$ws.Range("M122").Value = -6700

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 224
$ws.Range("I16").Value = 230.55556
$ws.Range("J16").Value = 165
$ws.Range("K16").Value = 230.55556
$ws.Range("L16").Value = 165
$ws.Range("M16").Value = -60.55556000000001
$ws.Range("N16").Value = -505
# Row 40
$ws.Range("H40").Value = 4558.7856
$ws.Range("I40").Value = 4558.7856
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4558.7856
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4422.7856
$ws.Range("N40").ClearContents()
# Row 100
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
